$d = $word.ActiveDocument

# Build the replacement OOXML for the document body content.
# This preserves all of the original (unchanged) run/paragraph attributes
# exactly, and adds:
#  1. <w:proofErr w:type="gramStart"/> / <w:proofErr w:type="gramEnd"/>
#     around the '" ."' sequence at the end of the first paragraph.
#  2. A paragraph break right after "...forma correcta." splitting the
#     single paragraph into two, moving the _GoBack bookmark into the
#     new second paragraph.
#  3. A brand-new second paragraph with the 14/02/2014 log entry, complete
#     with spelling/grammar proofErr markers matching the source diff.

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="003D1985" w:rsidRDefault="006A624A"><w:r><w:t>13/02/2014 2</w:t></w:r><w:r w:rsidR="00D00665"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00D00665"><w:t>hr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="009919EF"><w:t xml:space="preserve">. Modificadas las </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="009919EF"><w:t>Invalid</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="009919EF"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="009919EF"><w:t>Expression</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="009919EF"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="009919EF"><w:t>Exception</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="009919EF"><w:t xml:space="preserve"> y agregada aclaración que se debe usar </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="009919EF"><w:t>Typedef</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="009919EF"><w:t xml:space="preserve"> en las pruebas.</w:t></w:r><w:r><w:t>&#8221;</w:t></w:r><w:r w:rsidR="00A7169F"><w:t xml:space="preserve"> -&gt;</w:t></w:r><w:r><w:t>&#8221;</w:t></w:r><w:r w:rsidR="00A7169F"><w:t xml:space="preserve"> será tratado como </w:t></w:r><w:r><w:t>&#8220;</w:t></w:r><w:r w:rsidR="00A7169F"><w:t>.</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>&#8221;</w:t></w:r><w:r w:rsidR="00A7169F"><w:t xml:space="preserve"> .</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> Ahora puede leer varios atributos declarados en la misma línea de forma correcta.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">14/02/2014 3hr. Ya guarda los atributos </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>globales</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>,pero</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> todavía no los incluye en las heurísticas ni los muestra. No corre todavía con </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>Aeropuerto.c</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> .</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> Ahora lee números de </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>mas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> de 1 cifra. Ahora puede leer caracteres encerrados entre apostrofes. Ahora soporta funciones que devuelven estructuras.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$d.Content.InsertXML($xml)
